# daily auto push: 2026-01-27 07:38 UTC
# Inserts a new reading (2026/01/27 -> 15 o'clock) ahead of the existing
# 2026/12/29 block, shifting every later row down by one and appending the
# final trailing row, growing the sheet from A1:D765 to A1:D766.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 724
$endRow = 766

# Column A holds plain-text dates (e.g. "2026/01/27"); force text formatting
# first so Excel doesn't auto-coerce the strings into date serials.
$ws.Range("A" + $startRow + ":A" + $endRow).NumberFormat = "@"

$data = @(
    @("2026/01/27","火",15),
    @("2026/01/27","火",13),
    @("2026/01/27","火",16),
    @("2026/01/27","火",19),
    @("2026/12/29","火",23),
    @("2026/12/29","火",2),
    @("2026/12/29","火",5),
    @("2026/12/29","火",8),
    @("2026/12/29","火",13),
    @("2026/12/29","火",16),
    @("2026/12/30","水",22),
    @("2026/12/30","水",2),
    @("2026/12/30","水",6),
    @("2026/12/30","水",10),
    @("2026/12/30","水",12),
    @("2026/12/30","水",14),
    @("2026/12/31","木",22),
    @("2026/12/31","木",2),
    @("2026/12/31","木",5),
    @("2026/12/31","木",13),
    @("2026/12/31","木",16),
    @("2027/01/01","金",19),
    @("2027/01/01","金",1),
    @("2027/01/01","金",5),
    @("2027/01/01","金",8),
    @("2027/01/01","金",13),
    @("2027/01/01","金",16),
    @("2027/01/01","金",19),
    @("2027/01/02","土",22),
    @("2027/01/02","土",1),
    @("2027/01/02","土",4),
    @("2027/01/02","土",7),
    @("2027/01/02","土",13),
    @("2027/01/02","土",16),
    @("2027/01/02","土",19),
    @("2027/01/03","日",22),
    @("2027/01/03","日",2),
    @("2027/01/03","日",4),
    @("2027/01/03","日",7),
    @("2027/01/03","日",13),
    @("2027/01/04","月",22),
    @("2027/01/04","月",2),
    @("2027/01/05","火",7)
)

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = 201
    $r = $r + 1
}
